$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 361
$ws.Range("D2").Value = 39
$ws.Range("B5").Value = 0.9025
$ws.Range("D5").Value = 0.0975
